$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Mã phân công 1 -> 1011, advisor "Nguyễn Văn A" -> "Nguyễn Đức Trung"
$ws.Range("A2").Value = 1011
$ws.Range("B2").Value = "D21CQCN01-N"
$ws.Range("C2").Value = "2023-2024"
$ws.Range("D2").Value = "Nguyễn Đức Trung"
$ws.Range("E2").Value = 1

# Row 3: Mã phân công 2 -> 1012, class D21CQCN02-N -> D21CQCN01-N,
# year 2023-2024 -> 2022-2023, advisor "Nguyễn Văn A" -> "Phan Văn Anh"
$ws.Range("A3").Value = 1012
$ws.Range("B3").Value = "D21CQCN01-N"
$ws.Range("C3").Value = "2022-2023"
$ws.Range("D3").Value = "Phan Văn Anh"
$ws.Range("E3").Value = 1

# Row 4: Mã phân công 5 -> 1013, class D21CQCN01-N -> D21CQCN02-N,
# advisor "Trần Minh Hiếu" -> "Nguyễn Đức Trung"
$ws.Range("A4").Value = 1013
$ws.Range("B4").Value = "D21CQCN02-N"
$ws.Range("C4").Value = "2023-2024"
$ws.Range("D4").Value = "Nguyễn Đức Trung"
$ws.Range("E4").Value = 1

# Row 5: Mã phân công 8 -> 1021, class D21CQCN01-N -> D21CQCN02-N,
# year 2022-2023 (unchanged text), advisor "Trần Minh Hiếu" -> "Nguyễn Đức Trung",
# Trạng thái hiển thị 0 -> 1
$ws.Range("A5").Value = 1021
$ws.Range("B5").Value = "D21CQCN02-N"
$ws.Range("C5").Value = "2022-2023"
$ws.Range("D5").Value = "Nguyễn Đức Trung"
$ws.Range("E5").Value = 1
